$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Range('D2').Value = '26.799.57'
$ws.Range('E2').Value = '  -1.58%  '
$ws.Range('D3').Value = '1.871.31'
$ws.Range('E3').Value = '  -1.78%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '300.91'
$ws.Range('E6').Value = '  -0.06%  '
$ws.Range('D7').Value = '0.5344'
$ws.Range('E7').Value = '  +1.73%  '
$ws.Range('D8').Value = '0.3744'
$ws.Range('E8').Value = '  -1.68%  '
$ws.Range('D9').Value = '0.07185'
$ws.Range('E9').Value = '  -1.49%  '
$ws.Range('D10').Value = '21.63'
$ws.Range('E10').Value = '  -0.40%  '
$ws.Range('D11').Value = '0.8890'
$ws.Range('E11').Value = '  -1.54%  '
$ws.Range('D12').Value = '0.08165'
$ws.Range('E12').Value = '  -0.64%  '
$ws.Range('D13').Value = '1.879.57'
$ws.Range('E13').Value = '  +27.32%  '
$ws.Range('D14').Value = '93.00'
$ws.Range('E14').Value = '  -3.44%  '
$ws.Range('D15').Value = '5.310'
$ws.Range('E15').Value = '  -0.79%  '
$ws.Range('D16').Value = '1.001'
$ws.Range('E16').Value = '  -0.10%  '
$ws.Range('D17').Value = '14.84'
$ws.Range('E17').Value = '  +0.57%  '
$ws.Range('D18').Value = '0.000008522'
$ws.Range('E18').Value = '  -1.57%  '
$ws.Range('E19').Value = '  -0.07%  '
$ws.Range('D20').Value = '26.835.81'
$ws.Range('E20').Value = '  -1.59%  '
$ws.Range('D21').Value = '4.982'
$ws.Range('E21').Value = '  -2.51%  '
$ws.Range('D22').Value = '10.63'
$ws.Range('E22').Value = '  -1.65%  '
$ws.Range('D23').Value = '6.380'
$ws.Range('E23').Value = '  -1.90%  '
$ws.Range('D24').Value = '2.310'
$ws.Range('E24').Value = '  -1.17%  '
$ws.Range('D25').Value = '146.00'
$ws.Range('E25').Value = '  -2.84%  '
$ws.Range('D26').Value = '1.731'
$ws.Range('E26').Value = '  -0.68%  '
$ws.Range('E27').Value = '  -1.19%  '
$ws.Range('D28').Value = '113.86'
$ws.Range('E28').Value = '  -2.38%  '
$ws.Range('D29').Value = '4.723'
$ws.Range('E29').Value = '  -2.42%  '
$ws.Range('D30').Value = '4.625'
$ws.Range('E30').Value = '  -4.78%  '
$ws.Range('E31').Value = '  -1.05%  '
$ws.Range('D32').Value = '0.8045'
$ws.Range('E32').Value = '  -2.98%  '
$ws.Range('D33').Value = '0.05026'
$ws.Range('E33').Value = '  -0.57%  '
$ws.Range('E34').Value = '  -4.16%  '
$ws.Range('D35').Value = '2.941'
$ws.Range('E35').Value = '  -1.33%  '
$ws.Range('E36').Value = '  +5.89%  '
$ws.Range('D37').Value = '2.693'
$ws.Range('E37').Value = '  -1.18%  '
$ws.Range('D38').Value = '3.189'
$ws.Range('E38').Value = '  -4.93%  '
$ws.Range('D39').Value = '0.01953'
$ws.Range('E39').Value = '  -2.61%  '
$ws.Range('E40').Value = '  -1.43%  '
$ws.Range('D41').Value = '6.537'
$ws.Range('E41').Value = '  -1.17%  '
$ws.Range('D42').Value = '0.5204'
$ws.Range('E42').Value = '  +5.93%  '
$ws.Range('D43').Value = '8.775'
$ws.Range('E43').Value = '  -4.74%  '
$ws.Range('D44').Value = '114.58'
$ws.Range('E44').Value = '  -2.15%  '
$ws.Range('D45').Value = '0.1493'
$ws.Range('E45').Value = '  -1.72%  '
$ws.Range('D47').Value = '1.644'
$ws.Range('E47').Value = '  +0.14%  '
$ws.Range('D48').Value = '9.965'
$ws.Range('E48').Value = '  -2.75%  '
$ws.Range('E49').Value = '  -3.27%  '
$ws.Range('D50').Value = '0.06056'
$ws.Range('E50').Value = '  -0.11%  '
$ws.Range('D51').Value = '62.16'
$ws.Range('E51').Value = '  -3.51%  '
